$d = $word.ActiveDocument

# --- Merge split runs in Title / Author / Abstract paragraphs into single runs ---
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Title") { $para = $p; break }
}
$full = $para.Range
$target = $d.Range($full.Start, $full.End - 1)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Answers: Vector addition and scalar multiplication</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($frag)

$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Author") { $para = $p; break }
}
$full = $para.Range
$target = $d.Range($full.Start, $full.End - 1)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Renee Knapp, Kin Wang Pang</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($frag)

$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Abstract") { $para = $p; break }
}
$full = $para.Range
$target = $d.Range($full.Start, $full.End - 1)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Answers to questions relating to the guide on vector addition and scalar multiplication.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($frag)

# --- Fix m:dPr child-element order (begChr, sepChr, endChr, grow) in math zones ---
$om = $d.OMaths.Item(11)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>11</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>10</m:t></m:r></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(13)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:t>12</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>19</m:t></m:r></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(15)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>7</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>9</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>14</m:t></m:r></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(17)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>4</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>7</m:t></m:r><m:r><m:t>y</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(18)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>a</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>7</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>y</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>z</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(20)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(22)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>u</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:d><m:r><m:t>5</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>j</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r></m:e></m:d><m:r><m:t>6</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>k</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>15</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>j</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>18</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>k</m:t></m:r></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(23)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>6</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>v</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>18</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>42</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(24)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>v</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>u</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>27</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>10</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(25)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>w</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>4</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>u</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="b" /></m:rPr><m:t>v</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>32</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(32)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:acc><m:accPr><m:chr m:val="⃗" /></m:accPr><m:e><m:r><m:t>A</m:t></m:r><m:r><m:t>B</m:t></m:r></m:e></m:acc><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>7</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:mr></m:m></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>2</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(33)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:acc><m:accPr><m:chr m:val="⃗" /></m:accPr><m:e><m:r><m:t>A</m:t></m:r><m:r><m:t>B</m:t></m:r></m:e></m:acc><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>6</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(34)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:acc><m:accPr><m:chr m:val="⃗" /></m:accPr><m:e><m:r><m:t>A</m:t></m:r><m:r><m:t>C</m:t></m:r></m:e></m:acc><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(35)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:acc><m:accPr><m:chr m:val="⃗" /></m:accPr><m:e><m:r><m:t>A</m:t></m:r><m:r><m:t>B</m:t></m:r></m:e></m:acc><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:acc><m:accPr><m:chr m:val="⃗" /></m:accPr><m:e><m:r><m:t>A</m:t></m:r><m:r><m:t>C</m:t></m:r></m:e></m:acc><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>6</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>10</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>5</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(37)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:acc><m:accPr><m:chr m:val="⃗" /></m:accPr><m:e><m:r><m:t>C</m:t></m:r><m:r><m:t>B</m:t></m:r></m:e></m:acc><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>6</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>0</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>11</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>7</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr></m:m></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>6</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>10</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>5</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(39)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>5</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>9</m:t></m:r></m:e></m:mr></m:m></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:sSub><m:e><m:r><m:t>a</m:t></m:r></m:e><m:sub><m:r><m:t>1</m:t></m:r></m:sub></m:sSub></m:e></m:mr><m:mr><m:e><m:sSub><m:e><m:r><m:t>a</m:t></m:r></m:e><m:sub><m:r><m:t>2</m:t></m:r></m:sub></m:sSub></m:e></m:mr><m:mr><m:e><m:sSub><m:e><m:r><m:t>a</m:t></m:r></m:e><m:sub><m:r><m:t>3</m:t></m:r></m:sub></m:sSub></m:e></m:mr></m:m></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>6</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>7</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(40)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:t>A</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>,</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>,</m:t></m:r><m:r><m:t>11</m:t></m:r></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

$om = $d.OMaths.Item(47)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMath><m:r><m:t>2</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>2</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>5</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>z</m:t></m:r></m:e></m:mr></m:m></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>4</m:t></m:r></m:e></m:mr></m:m></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:m><m:mPr><m:baseJc m:val="center" /><m:plcHide m:val="on" /><m:mcs><m:mc><m:mcPr><m:mcJc m:val="center" /><m:count m:val="1" /></m:mcPr></m:mc></m:mcs></m:mPr><m:mr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>y</m:t></m:r></m:e></m:mr><m:mr><m:e><m:r><m:t>0</m:t></m:r></m:e></m:mr></m:m></m:e></m:d></m:oMath></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$om.Range.InsertXML($frag)

